$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1744
$ws.Range("I19").Value = 1498.1111
$ws.Range("J19").Value = 1928.4166
$ws.Range("K19").Value = 1498.1111
$ws.Range("L19").Value = 1928.4166
$ws.Range("M19").Value = -1323.1111
$ws.Range("N19").Value = -2278.4166

# Hunk 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 90123.38
$ws.Range("I62").Value = 113420.4
$ws.Range("J62").Value = 12466.667
$ws.Range("K62").Value = 113420.4
$ws.Range("L62").Value = 12466.667
$ws.Range("M62").Value = -112796.4
$ws.Range("N62").Value = -13714.667

# Hunk 2: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 90123.38
$ws.Range("I65").Value = 113420.4
$ws.Range("J65").Value = 12466.667
$ws.Range("K65").Value = 567102
$ws.Range("L65").Value = 62333.335
$ws.Range("M65").Value = -563982
$ws.Range("N65").Value = -68573.33499999999

# Hunk 3: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 951.5625
$ws.Range("I98").Value = 1032.6923
$ws.Range("K98").Value = 1032.6923
$ws.Range("M98").Value = 465.3077000000001

# Hunk 4: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 951.5625
$ws.Range("I122").Value = 1032.6923
$ws.Range("K122").Value = 3098.0769
$ws.Range("M122").Value = -648.0769

# Hunk 5: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3491.44
$ws.Range("I138").Value = 904.5333000000001
$ws.Range("J138").Value = 3947.953
$ws.Range("K138").Value = 2713.5999
$ws.Range("L138").Value = 11843.859
$ws.Range("M138").Value = 2426.4001
$ws.Range("N138").Value = -22123.859

# Hunk 6: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("K3").Value = 5
$ws.Range("M3").Value = 110

# Hunk 7: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1476.4242
$ws.Range("I61").Value = 1319.84
$ws.Range("J61").Value = 1965.75
$ws.Range("K61").Value = 1319.84
$ws.Range("L61").Value = 1965.75
$ws.Range("M61").Value = -1107.84
$ws.Range("N61").Value = -2389.75

# Hunk 8: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1494.1
$ws.Range("I122").Value = 1436.3334
$ws.Range("J122").Value = 2014
$ws.Range("K122").Value = 4309.0002
$ws.Range("L122").Value = 6042
$ws.Range("M122").Value = -1859.0002
$ws.Range("N122").Value = -10942

# Hunk 9: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1755.4445
$ws.Range("I132").Value = 1249.8462
$ws.Range("J132").Value = 2224.9285
$ws.Range("K132").Value = 3749.5386
$ws.Range("L132").Value = 6674.7855
$ws.Range("M132").Value = -1219.5386
$ws.Range("N132").Value = -11734.7855

# Hunk 10: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1476.4242
$ws.Range("I136").Value = 1319.84
$ws.Range("J136").Value = 1965.75
$ws.Range("K136").Value = 3959.52
$ws.Range("L136").Value = 5897.25
$ws.Range("M136").Value = -1409.52
$ws.Range("N136").Value = -10997.25

# Hunk 11: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27784.805
$ws.Range("I134").Value = 3616.4348
$ws.Range("J134").Value = 58666.61
$ws.Range("K134").Value = 10849.3044
$ws.Range("L134").Value = 175999.83
$ws.Range("M134").Value = -8314.304400000001
$ws.Range("N134").Value = -181069.83

# Hunk 12: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1887
$ws.Range("N3").ClearContents()

# Hunk 13: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1667387.1
$ws.Range("I122").Value = 1667387.1
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5002161.300000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4999711.300000001
$ws.Range("N122").ClearContents()

# Hunk 14: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2833.84
$ws.Range("I132").Value = 2002.8334
$ws.Range("J132").Value = 4970.7144
$ws.Range("K132").Value = 6008.5002
$ws.Range("L132").Value = 14912.1432
$ws.Range("M132").Value = -3478.5002
$ws.Range("N132").Value = -19972.1432

# Hunk 15: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3825
$ws.Range("I134").Value = 3320.4546
$ws.Range("J134").Value = 4750
$ws.Range("K134").Value = 9961.363799999999
$ws.Range("L134").Value = 14250
$ws.Range("M134").Value = -7426.363799999999
$ws.Range("N134").Value = -19320

# Hunk 16: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2049.6155
$ws.Range("I4").Value = 233.33333
$ws.Range("J4").Value = 2286.5217
$ws.Range("K4").Value = 699.99999
$ws.Range("L4").Value = 6859.5651
$ws.Range("M4").Value = -587.99999
$ws.Range("N4").Value = -7083.5651

# Hunk 17: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 633.3333
$ws.Range("I68").Value = 633.3333
$ws.Range("K68").Value = 1899.9999
$ws.Range("M68").Value = -1088.9999

# Hunk 18: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 633.3333
$ws.Range("I71").Value = 633.3333
$ws.Range("K71").Value = 5699.9997
$ws.Range("M71").Value = -1643.9997

# Hunk 19: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 371227.16
$ws.Range("J122").Value = 910396.25
$ws.Range("L122").Value = 8193566.25
$ws.Range("N122").Value = -8198466.25

# Hunk 20: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4576.0527
$ws.Range("I70").Value = 4076.8462
$ws.Range("J70").Value = 5657.6665
$ws.Range("K70").Value = 4076.8462
$ws.Range("L70").Value = 5657.6665
$ws.Range("M70").Value = -3806.8462
$ws.Range("N70").Value = -6197.6665

# Hunk 21: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4576.0527
$ws.Range("I73").Value = 4076.8462
$ws.Range("J73").Value = 5657.6665
$ws.Range("K73").Value = 4076.8462
$ws.Range("L73").Value = 5657.6665
$ws.Range("M73").Value = -3140.8462
$ws.Range("N73").Value = -7529.6665

# Hunk 22: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1437.1364
$ws.Range("I102").Value = 1178.7858
$ws.Range("J102").Value = 1889.25
$ws.Range("K102").Value = 1178.7858
$ws.Range("L102").Value = 1889.25
$ws.Range("M102").Value = 443.2141999999999
$ws.Range("N102").Value = -5133.25

# Hunk 23: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2350.0588
$ws.Range("I122").Value = 1733.2858
$ws.Range("K122").Value = 5199.857400000001
$ws.Range("M122").Value = -2749.857400000001

# Hunk 24: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2378.9
$ws.Range("I40").Value = 2345.7646
$ws.Range("J40").Value = 2566.6667
$ws.Range("K40").Value = 2345.7646
$ws.Range("L40").Value = 2566.6667
$ws.Range("M40").Value = -2209.7646
$ws.Range("N40").Value = -2838.6667

# Hunk 25: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3897.3103
$ws.Range("I122").Value = 5674.815
$ws.Range("J122").Value = 2349.1614
$ws.Range("K122").Value = 17024.445
$ws.Range("L122").Value = 7047.4842
$ws.Range("M122").Value = -14574.445
$ws.Range("N122").Value = -11947.4842

# Hunk 26: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2550.7144
$ws.Range("I132").Value = 1996.2858
$ws.Range("J132").Value = 4214
$ws.Range("K132").Value = 5988.857400000001
$ws.Range("L132").Value = 12642
$ws.Range("M132").Value = -3458.857400000001
$ws.Range("N132").Value = -17702

# Hunk 27: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1153.375
$ws.Range("I122").Value = 1113.2273
$ws.Range("K122").Value = 3339.6819
$ws.Range("M122").Value = -889.6819

# Hunk 28: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1710.275
$ws.Range("I132").Value = 1405.0968
$ws.Range("J132").Value = 2761.4443
$ws.Range("K132").Value = 4215.2904
$ws.Range("L132").Value = 8284.332900000001
$ws.Range("M132").Value = -1685.2904
$ws.Range("N132").Value = -13344.3329
